# Update "Corr/total marks" on the concise marksheet:
#   - Marking row (B11): correct-answer count 3 -> 5
#   - Total row (B12): total marks 72 -> 120
#   - Total row (E12): "obtained/max" label 71/84 -> 120/140
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
